# Weekly update: insert a new daily-price record as the new row 48
# (Vega Central Mapocho de Santiago - Arveja Verde), pushing the existing
# rows 48-70 down to 49-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 48; this shifts rows 48:70 down to 49:71
# and extends the used range to A1:R71, exactly like Excel's native
# "Insert Sheet Rows" command.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new record's data.
$ws.Cells.Item(48, 1).Value = 9
$ws.Cells.Item(48, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(48, 3).Value = "Metropolitana"
$ws.Cells.Item(48, 4).Value = 44489
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = 100112022
$ws.Cells.Item(48, 7).Value = "Arveja Verde"
$ws.Cells.Item(48, 8).Value = "Perfection"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 16
$ws.Cells.Item(48, 11).Value = 24000
$ws.Cells.Item(48, 12).Value = 25000
$ws.Cells.Item(48, 13).Value = 24500
$ws.Cells.Item(48, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(48, 16).Value = 980
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
